# Reorder the worksheet tabs: "review_info" moves in front of "hotel_info"
# (it becomes the first/active sheet, "hotel_info" becomes the second).
$wb = $excel.ActiveWorkbook

$reviewInfo = $wb.Worksheets.Item("review_info")
$reviewInfo.Move($wb.Worksheets.Item(1))

# Add a new "State" column to the hotel_info sheet, inserted right after
# "Hotel_Name" and before "City", with value "Louisiana" for the existing
# hotel record.
$hotelInfo = $wb.Worksheets.Item("hotel_info")
$hotelInfo.Range("C1").EntireColumn.Insert()
$hotelInfo.Range("C1").Value = "State"
$hotelInfo.Range("C2").Value = "Louisiana"
